# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (formatted like the existing "2021-Q4"
# sheet) between "2021-Q4" and "总计", and adds a corresponding summary row
# to "总计".

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" sheet right after "2021-Q4" ---------------
# NOTE: worksheet object references in this runtime are position-based, not
# identity-based, so any sheet reference obtained *before* inserting a new
# sheet may silently re-point to a different sheet afterwards. Grab
# "总计" fresh (by name) only after the insert is complete, below.
$q1Sheet = $wb.Worksheets.Add($null, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# Copy header-row (B1:H1) and index-cell (A2) formatting from the "2021-Q4"
# sheet so the new sheet matches the existing visual style (bold header,
# bordered A2 index cell) without duplicating style entries.
$q4Sheet.Range("B1:H1").Copy()
$q1Sheet.Range("B1").PasteSpecial(-4122)

$q4Sheet.Range("A2").Copy()
$q1Sheet.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Header text
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Row 2 data - A2 is the numeric row index (0), B2:G2 are stored as plain
# text (matching the source data's inlineStr cells), H2 is numeric.
$q1Sheet.Range("A2").Value = 0

$q1Sheet.Range("B2:G2").NumberFormat = "@"
$q1Sheet.Range("B2").Value = "968013"
$q1Sheet.Range("C2").Value = "施罗德亚洲高息股债基金M"
$q1Sheet.Range("D2").Value = "297.64"
$q1Sheet.Range("E2").Value = "57.54"
$q1Sheet.Range("F2").Value = "1.49"
$q1Sheet.Range("G2").Value = "4.4348"
$q1Sheet.Range("B2:G2").ClearFormats()

$q1Sheet.Range("H2").Value = 4

# --- 2. Update "总计": push the existing 2021-Q4 row to row 3, write a new
#        2022-Q1 row into row 2 ------------------------------------------
# Fetch this sheet reference now (by name), after the insert above, so it
# resolves to the correct (shifted) worksheet position.
$totalSheet = $wb.Worksheets.Item("总计")

# Give row 3's index cell (A3) the same style as A2 (bold/bordered) before
# moving the old row 2 values down, mirroring the existing A2 formatting.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the existing 2021-Q4 summary row down to row 3.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 3.06

# Write the new 2022-Q1 summary row into row 2 (A2 already carries the
# correct style from the original sheet; B2:D2 stay unstyled).
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 4.43

Write-Host "Done: 2022-Q1 sheet added; summary sheet updated."
